{"js": "// Insert a new list item \"challenge part----BMI calculator---\" (bold, red)\n// right after the \"Training and exercise----optional\" bullet, matching the\n// formatting (ListParagraph style, same numbering list) of its neighbours.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst anchorText = \"Training and exercise----optional\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not locate anchor paragraph containing \"' + anchorText + '\"');\n}\n\n// insertParagraph after the anchor inherits the anchor's paragraph/run\n// formatting (ListParagraph style + numbering + bold/red run), which is\n// exactly the formatting the new bullet needs.\nconst newPara = anchor.insertParagraph(\n  \"challenge part----BMI calculator---\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert a new list item \"challenge part----BMI calculator---\" (bold, red)\n# right after the \"Training and exercise----optional\" bullet, matching the\n# formatting (ListParagraph style, same numbering list) of its neighbours.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Training and exercise----optional\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $target = $p\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate anchor paragraph containing '$anchorText'\"\n}\n\n# InsertParagraphAfter() creates a new empty paragraph that inherits the\n# anchor's paragraph/run formatting (ListParagraph style + numbering +\n# bold/red run) - exactly the formatting the new bullet needs.\n$target.Range.InsertParagraphAfter()\n\n$newPara = $target.Next()\n$newPara.Range.Text = \"challenge part----BMI calculator---\"\n"}
